# feat: change PFAS selection criteria
# PFPeA and 6:2 FTSA no longer meet the PFAS selection criteria, so their
# rows are removed from every type/grp block, and the sumPFAS (min/median/max)
# rows are recomputed without their contribution.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the sumPFAS rows (min/median/max) with the recomputed values
#    BEFORE removing any rows, while row numbers still match the original layout.
$ws.Range("D2").Value = 4.083778437291977
$ws.Range("E2").Value = 5.489575741116197
$ws.Range("F2").Value = 9.162230259008391

$ws.Range("D16").Value = 11.07954611693446
$ws.Range("E16").Value = 15.03551176672183
$ws.Range("F16").Value = 17.70102230973641

$ws.Range("D30").Value = 14.11360378205541
$ws.Range("E30").Value = 21.16570361356836
$ws.Range("F30").Value = 25.7597628633788

$ws.Range("D44").Value = 14.32840609082056
$ws.Range("E44").Value = 17.1871208667622
$ws.Range("F44").Value = 25.79997209101526

# 2) Delete the PFPeA and "6:2 FTSA" rows for each group, from the bottom
#    up so earlier row numbers stay valid as later rows are removed.
$rowsToDelete = @(56, 45, 42, 31, 28, 17, 14, 3)
foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete()
}
